# "Menu Test Cases added"
#
# The test-data sheet gets 5 new rows appended to the existing login test
# cases: one more login scenario (WrongPassLogin) plus four menu-navigation
# smoke tests (OpenHomePage / OpenContactPage / OpenSalesPage /
# OpenPurchasePage), each expected to "Pass".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - another login test case (wrong password this time).
$ws.Range("A4").Value = "WrongPassLogin"
$ws.Range("B4").Value = "admin"
$ws.Range("C4").Value = " "
$ws.Range("D4").Value = "Pass"

# Rows 5-8 - new menu navigation test cases. Column A values are entered in
# this order so the shared-string table lines up with how the workbook was
# authored (WrongPassLogin, OpenContactPage, OpenSalesPage, OpenPurchasePage,
# then OpenHomePage).
$ws.Range("A6").Value = "OpenContactPage"
$ws.Range("A7").Value = "OpenSalesPage"
$ws.Range("A8").Value = "OpenPurchasePage"
$ws.Range("A5").Value = "OpenHomePage"

$ws.Range("D5").Value = "Pass"
$ws.Range("D6").Value = "Pass"
$ws.Range("D7").Value = "Pass"
$ws.Range("D8").Value = "Pass"

# Column D cells in the existing rows have no explicit cell style (they just
# use the column's default formatting) - reset the new ones to match rather
# than leaving the style picked up implicitly from neighboring cells.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"

# Column A now holds longer strings (e.g. "OpenPurchasePage") so it no longer
# fits the old best-fit width; give it an explicit custom width sized to the
# new content instead.
$ws.Columns("A").ColumnWidth = 18.833333

# Move the active selection, matching where the user left off editing.
$null = $ws.Range("E14").Select()
